$d = $word.ActiveDocument
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>2 – Drunken battle with thug in village</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>IF won</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:tab/>
        <w:t>(3.38) // Won in battle with thug</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>ELSE</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:tab/>
        <w:t>\0015\</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>3 – Fight against old chieftain</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="708"/>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/>
          <w:color w:val="000000"/>
          <w:sz w:val="17"/>
          <w:szCs w:val="17"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/>
          <w:color w:val="000000"/>
          <w:sz w:val="17"/>
          <w:szCs w:val="17"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>IF lost</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="708" w:firstLine="708"/>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/>
          <w:color w:val="000000"/>
          <w:sz w:val="17"/>
          <w:szCs w:val="17"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/>
          <w:color w:val="000000"/>
          <w:sz w:val="17"/>
          <w:szCs w:val="17"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>\0022\ // Killed by chieftain</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="708"/>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/>
          <w:color w:val="000000"/>
          <w:sz w:val="17"/>
          <w:szCs w:val="17"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/>
          <w:color w:val="000000"/>
          <w:sz w:val="17"/>
          <w:szCs w:val="17"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>IF won</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="708" w:firstLine="708"/>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/>
          <w:color w:val="000000"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/>
          <w:color w:val="000000"/>
          <w:sz w:val="17"/>
          <w:szCs w:val="17"/>
          <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>(9.4) // Chieftain is killed.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>4 – Fight with artifact guard</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:tab/>
        <w:t>IF won</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:tab/>
        <w:t>(11.5) // Choise</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:tab/>
        <w:t>IF lost</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:tab/>
        <w:t>\0025\ // Killed by guard.</w:t>
      </w:r>
    </w:p>
'@
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertXML($xml)
